$d = $word.ActiveDocument

# Initial layout:
#   [1] "Iniciando Word"
#   [2] (empty)
#   [3] (empty, contains the _GoBack bookmark)
#
# Target layout:
#   [1]  "Iniciando Word"          (unchanged)
#   [2]  (empty)                   (unchanged)
#   [3]  (empty)                   (new)
#   [4]  (empty)                   (new)
#   [5]  (empty)                   (new)
#   [6]  "Tabla de contenido"      (new)
#   [7]  (empty)                   (new)
#   [8]  (empty, bookmark)         (unchanged paragraph, shifted down)
#   [9]  (empty)                   (new)
#   [10] (empty)                   (new)
#
# The paragraph that carries the _GoBack bookmark starts out as paragraph 3.
# Each time we insert a new paragraph immediately *before* it, that bookmark
# paragraph's own index moves one further down, so re-fetching
# Paragraphs.Item() right before every insert (instead of reusing a cached
# Range) keeps the indices reliable.

# Three new empty paragraphs before the bookmark paragraph (currently #3).
for ($i = 0; $i -lt 3; $i++) {
    $pBookmark = $d.Paragraphs.Item(3)
    $pBookmark.Range.InsertParagraphBefore()
}

# The bookmark paragraph is now #6. Insert one more empty paragraph before
# it (will become the "Tabla de contenido" paragraph) and fill it in.
$pBookmark = $d.Paragraphs.Item(6)
$pBookmark.Range.InsertParagraphBefore()
$pContents = $d.Paragraphs.Item(6)
$pContents.Range.Text = "Tabla de contenido"

# The bookmark paragraph is now #7. Insert one trailing empty paragraph
# before it.
$pBookmark = $d.Paragraphs.Item(7)
$pBookmark.Range.InsertParagraphBefore()

# The bookmark paragraph is now #8. Append two empty paragraphs after it.
$pBookmark = $d.Paragraphs.Item(8)
$pBookmark.Range.InsertParagraphAfter()
$pBookmark = $d.Paragraphs.Item(8)
$pBookmark.Range.InsertParagraphAfter()
